$wb = $excel.ActiveWorkbook

# --- Sheet 1: question_answers ---
# Column B holds the respondents chosen answer option, stored as TEXT ("1".."4").
$ws1 = $wb.Worksheets.Item("question_answers")

$answers = @{
    "2" = "1"
    "3" = "1"
    "4" = "2"
    "5" = "2"
    "6" = "2"
    "7" = "2"
    "10" = "1"
    "11" = "1"
    "12" = "3"
    "13" = "1"
    "14" = "4"
    "15" = "3"
    "16" = "4"
    "17" = "2"
    "19" = "3"
    "20" = "2"
    "21" = "4"
    "22" = "4"
    "24" = "3"
    "25" = "3"
    "26" = "3"
    "28" = "4"
    "30" = "3"
    "32" = "1"
    "34" = "1"
    "35" = "4"
    "38" = "1"
    "40" = "1"
    "41" = "3"
    "42" = "3"
    "43" = "4"
    "44" = "1"
    "45" = "3"
    "47" = "2"
    "48" = "2"
    "50" = "3"
    "51" = "2"
    "52" = "4"
    "53" = "4"
    "54" = "4"
    "55" = "1"
    "56" = "2"
    "57" = "4"
    "58" = "2"
    "60" = "4"
    "61" = "3"
    "62" = "3"
    "63" = "4"
    "64" = "1"
    "67" = "2"
    "69" = "2"
    "71" = "1"
    "72" = "4"
    "73" = "2"
    "74" = "3"
    "75" = "4"
    "77" = "2"
    "78" = "2"
    "79" = "2"
    "82" = "4"
    "83" = "4"
    "84" = "1"
    "85" = "2"
    "87" = "4"
    "89" = "2"
    "90" = "4"
    "91" = "1"
    "92" = "4"
    "93" = "3"
    "94" = "1"
    "95" = "3"
    "97" = "2"
    "98" = "4"
    "99" = "2"
    "100" = "1"
    "101" = "3"
    "105" = "3"
    "107" = "4"
    "108" = "3"
    "109" = "3"
    "110" = "4"
    "111" = "4"
}

foreach ($row in $answers.Keys) {
    $cell = $ws1.Range("B$row")
    $cell.NumberFormat = "@"
    $cell.Value = $answers[$row]
}

# --- Sheet 2: outputs ---
# Re-orders the dysthymia_type_a/b rows up next to dysthymia_total, and refreshes
# the recomputed score counts in column B.
$ws2 = $wb.Worksheets.Item("outputs")

$labels = @{
    "21" = "dysthymia_type_a"
    "22" = "dysthymia_type_b"
    "23" = "dysthymia_total"
    "24" = "autistic_disorder_type_a"
    "25" = "autistic_disorder_type_b"
    "26" = "autistic_disorder_type_c"
    "27" = "autistic_disorder_total"
    "28" = "asperger_disorder_type_a"
    "29" = "asperger_disorder_type_b"
    "30" = "asperger_disorder_total"
    "31" = "social_phobia"
    "32" = "seperation_anxiety_disorder"
    "33" = "enuresis"
}

foreach ($row in $labels.Keys) {
    $ws2.Range("A$row").Value = $labels[$row]
}

$scores = @{
    "2" = 6
    "3" = 7
    "5" = 6
    "6" = 14
    "8" = 0
    "9" = 4
    "21" = 1
    "22" = 2
    "25" = 4
    "27" = 11
    "28" = 3
    "29" = 4
    "30" = 7
    "31" = 2
    "32" = 6
    "33" = 1
    "37" = 1
    "38" = 0
    "43" = 1
    "48" = 1
    "49" = 1
}

foreach ($row in $scores.Keys) {
    $ws2.Range("B$row").Value = $scores[$row]
}
